$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.403.86'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '3.169.74'
$ws.Range('E3').Value = '  -4.24%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '''586.91'
$ws.Range('E5').Value = '  -2.70%  '
$ws.Range('D6').Value = '''135.03'
$ws.Range('E6').Value = '  -5.42%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.171.08'
$ws.Range('E8').Value = '  -4.16%  '
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('E10').Value = '  -5.43%  '
$ws.Range('D11').Value = '''5.25'
$ws.Range('E11').Value = '  -5.19%  '
$ws.Range('E12').Value = '  -4.43%  '
$ws.Range('E13').Value = '  -5.64%  '
$ws.Range('D14').Value = '''33.12'
$ws.Range('E14').Value = '  -4.73%  '
$ws.Range('D15').Value = '3.677.28'
$ws.Range('E15').Value = '  -4.62%  '
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').Value = '3.174.42'
$ws.Range('E17').Value = '  -4.12%  '
$ws.Range('D18').Value = '62.392.62'
$ws.Range('E18').Value = '  -2.41%  '
$ws.Range('E19').Value = '  -5.56%  '
$ws.Range('D20').Value = '''454.75'
$ws.Range('E20').Value = '  -5.52%  '
$ws.Range('D21').Value = '''13.87'
$ws.Range('E21').Value = '  -2.33%  '
$ws.Range('E22').Value = '  -4.53%  '
$ws.Range('D23').Value = '''7.60'
$ws.Range('E23').Value = '  -5.30%  '
$ws.Range('D24').Value = '''13.36'
$ws.Range('E24').Value = '  -0.89%  '
$ws.Range('D25').Value = '''83.15'
$ws.Range('E25').Value = '  -2.17%  '
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('E29').Value = '  -5.76%  '
$ws.Range('D30').Value = '''7.76'
$ws.Range('E30').Value = '  -4.97%  '
$ws.Range('E31').Value = '  -6.90%  '
$ws.Range('D32').Value = '''27.30'
$ws.Range('E32').Value = '  -7.09%  '
$ws.Range('E33').Value = '  -1.70%  '
$ws.Range('E34').Value = '  -6.42%  '
$ws.Range('E35').Value = '  -6.21%  '
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('E37').Value = '  -3.19%  '
$ws.Range('D38').Value = '0.0₃0691'
$ws.Range('E38').Value = '  -8.09%  '
$ws.Range('E39').Value = '  -4.78%  '
$ws.Range('D40').Value = '''410.92'
$ws.Range('E40').Value = '  -4.34%  '
$ws.Range('D41').Value = '2.877.00'
$ws.Range('E41').Value = '  -5.48%  '
$ws.Range('E42').Value = '  -3.28%  '
$ws.Range('D43').Value = '''7.98'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('E45').Value = '  +1.39%  '
$ws.Range('E46').Value = '  -6.47%  '
$ws.Range('E47').Value = '  -3.04%  '
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = '''124.60'
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('D50').Value = '''25.27'
$ws.Range('E50').Value = '  -4.32%  '
$ws.Range('D51').Value = '''0.111'
$ws.Range('E51').Value = '  -3.77%  '
